$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LINE_TRIALS_URL")

# --- Clear the old sample row (row 2), keeping the date formatting on D2:E2 ---
$ws.Range("A2:C2").ClearContents()
$ws.Range("D2:E2").ClearContents()
$ws.Range("F2:G2").ClearContents()

# --- Add the new trial rows (3-6) ---
$ws.Range("A3").Value2 = "PT NUSA SOLAR INDONESIA"
$ws.Range("B3").Value2 = "SOLARCELL"
$ws.Range("C3").Value2 = "In Progress"
$ws.Range("D3").Value2 = 45782
$ws.Range("E3").Value2 = "TBD"
$ws.Range("F3").Value2 = 22400

$ws.Range("A4").Value2 = "Dongguan CSG Solar Galss Co. L"
$ws.Range("B4").Value2 = "TEMPERED GLASS"
$ws.Range("C4").Value2 = "In Progress"
$ws.Range("D4").Value2 = 45782
$ws.Range("E4").Value2 = "TBD"
$ws.Range("F4").Value2 = 22400

$ws.Range("A5").Value2 = "Zhejiang Jiaxing Taihe New Energy"
$ws.Range("B5").Value2 = "AL FRAME LONG & SHORT"
$ws.Range("C5").Value2 = "In Progress"
$ws.Range("D5").Value2 = 45782
$ws.Range("E5").Value2 = "TBD"
$ws.Range("F5").Value2 = 22400

$ws.Range("A6").Value2 = "Suzhou UKT New Energy Technology"
$ws.Range("B6").Value2 = "JB SPLIT 1500V"
$ws.Range("C6").Value2 = "In Progress"
$ws.Range("D6").Value2 = 45782
$ws.Range("E6").Value2 = "TBD"
$ws.Range("F6").Value2 = 22400

# --- Copy the date number format from D2 onto the new start-date cells (D3:D6) ---
$ws.Range("D2").Copy()
$ws.Range("D3:D6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Loosen the data validation on the pasted END_DATE cells (E3:E19) ---
$ws.Range("E3:E19").Validation.Delete()
$ws.Range("E3:E19").Validation.Add(0)  # xlValidateInputOnly (no restriction)

# --- Restore the selection to the newly entered row, like after data entry ---
$ws.Range("A2:G2").Select()
